$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used in the existing palette.
$orange = 49407   # RGB(255,192,0) == FFC000, the existing "current stage" highlight
$white  = 2        # xlThemeColorLight1 -> OOXML theme index 0 (Background 1)

# --- Row 3 ("Cambiar las validaciones..."): progresses from "Asignado" to "En proceso" ---
$ws.Range("C3").Interior.ThemeColor = $white
$ws.Range("D3").Interior.Color = $orange

# --- Row 5 ("Mejorar los mensajes..."): progresses from "Asignado" to "En proceso" ---
$ws.Range("C5").Interior.ThemeColor = $white
$ws.Range("D5").Interior.Color = $orange

# --- Row 8 ("Revisar las excepciones..."): progresses from "Asignado" to "En proceso" ---
$ws.Range("C8").Interior.ThemeColor = $white
$ws.Range("D8").Interior.Color = $orange

# --- Row 9 ("Revisar y preguntar..."): progresses from "Asignado" to "En proceso" ---
$ws.Range("C9").Interior.ThemeColor = $white
$ws.Range("D9").Interior.Color = $orange
# Clear E9 to the same "passed" white fill, matching C3's already-established format
# (copy/paste the format instead of re-assigning ThemeColor, to avoid spawning a
# redundant intermediate fill entry when starting from an unfilled cell).
$ws.Range("C3").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 10 ("Funcionalidad generar reporte..."): assigned to Cristian Fernández,
#     and jumps straight to "Finalizado" ---
$ws.Range("B10").Value = "Cristian Fernández"
$ws.Range("E10").Interior.Color = $orange

# Restore the active selection to match the saved view state.
$ws.Range("F9").Select() | Out-Null
